$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 113
$ws.Range("I8").Value = 113
$ws.Range("K8").Value = 339
$ws.Range("M8").Value = -200

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 88
$ws.Range("I42").Value = 88
$ws.Range("K42").Value = 264
$ws.Range("M42").Value = -34

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1310.2972
$ws.Range("I112").Value = 1027.6666
$ws.Range("J112").Value = 1401.1428
$ws.Range("K112").Value = 3082.9998
$ws.Range("L112").Value = 4203.428400000001
$ws.Range("M112").Value = -1974.9998
$ws.Range("N112").Value = -6419.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1857.9259
$ws.Range("I127").Value = 588.05884
$ws.Range("J127").Value = 2441.3784
$ws.Range("K127").Value = 1764.17652
$ws.Range("L127").Value = 7324.135200000001
$ws.Range("M127").Value = 3195.82348
$ws.Range("N127").Value = -17244.1352

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1048.875
$ws.Range("I137").Value = 1005.8571
$ws.Range("K137").Value = 3017.5713
$ws.Range("M137").Value = -467.5712999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1251.0555
$ws.Range("I61").Value = 840.7
$ws.Range("J61").Value = 1764
$ws.Range("K61").Value = 840.7
$ws.Range("L61").Value = 1764
$ws.Range("M61").Value = -628.7
$ws.Range("N61").Value = -2188

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2551.7144
$ws.Range("I102").Value = 2042.5
$ws.Range("J102").Value = 3824.75
$ws.Range("K102").Value = 2042.5
$ws.Range("L102").Value = 3824.75
$ws.Range("M102").Value = -420.5
$ws.Range("N102").Value = -7068.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 682.619
$ws.Range("I122").Value = 662.5
$ws.Range("K122").Value = 1987.5
$ws.Range("M122").Value = 462.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1251.0555
$ws.Range("I136").Value = 840.7
$ws.Range("J136").Value = 1764
$ws.Range("K136").Value = 2522.1
$ws.Range("L136").Value = 5292
$ws.Range("M136").Value = 27.89999999999964
$ws.Range("N136").Value = -10392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 577
$ws.Range("J12").Value = 980
$ws.Range("L12").Value = 980
$ws.Range("N12").Value = -1320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4226.625
$ws.Range("I62").Value = 4421.8184
$ws.Range("J62").Value = 3797.2
$ws.Range("K62").Value = 4421.8184
$ws.Range("L62").Value = 3797.2
$ws.Range("M62").Value = -3797.8184
$ws.Range("N62").Value = -5045.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 4226.625
$ws.Range("I65").Value = 4421.8184
$ws.Range("J65").Value = 3797.2
$ws.Range("K65").Value = 22109.092
$ws.Range("L65").Value = 18986
$ws.Range("M65").Value = -18989.092
$ws.Range("N65").Value = -25226

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3272.2222
$ws.Range("I99").Value = 3150
$ws.Range("J99").Value = 3516.6667
$ws.Range("K99").Value = 3150
$ws.Range("L99").Value = 3516.6667
$ws.Range("M99").Value = -1652
$ws.Range("N99").Value = -6512.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3272.2222
$ws.Range("I126").Value = 3150
$ws.Range("J126").Value = 3516.6667
$ws.Range("K126").Value = 9450
$ws.Range("L126").Value = 10550.0001
$ws.Range("M126").Value = -6980
$ws.Range("N126").Value = -15490.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 888.6818
$ws.Range("I134").Value = 603.4
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 1810.2
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = 724.8000000000002
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 8650.5
$ws.Range("I44").Value = 325.75
$ws.Range("J44").Value = 25300
$ws.Range("K44").Value = 977.25
$ws.Range("L44").Value = 75900
$ws.Range("M44").Value = -579.25
$ws.Range("N44").Value = -76696

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 10795
$ws.Range("I141").Value = 11500
$ws.Range("K141").Value = 34500
$ws.Range("M141").Value = -29320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2211.111
$ws.Range("I102").Value = 1800
$ws.Range("K102").Value = 1800
$ws.Range("M102").Value = -178

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18521686
$ws.Range("I7").Value = 4125
$ws.Range("J7").Value = 55556810
$ws.Range("K7").Value = 4125
$ws.Range("L7").Value = 55556810
$ws.Range("M7").Value = -4013
$ws.Range("N7").Value = -55557034

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 18521686
$ws.Range("I126").Value = 4125
$ws.Range("J126").Value = 55556810
$ws.Range("K126").Value = 12375
$ws.Range("L126").Value = 166670430
$ws.Range("M126").Value = -9905
$ws.Range("N126").Value = -166675370

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4833.3335
$ws.Range("I96").Value = 4750
$ws.Range("K96").Value = 4750
$ws.Range("M96").Value = -3377

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 1028380
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 1028380
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 1028380
$ws.Range("N119").Value = -1038056

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 31250.25
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 31250.25
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 31250.25
$ws.Range("N120").Value = -40926.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 24000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 24000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 24000
$ws.Range("N121").Value = -27494

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1384.8334
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 1461
$ws.Range("K122").Value = 3012
$ws.Range("L122").Value = 4383
$ws.Range("M122").Value = -562
$ws.Range("N122").Value = -9283

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 27074.875
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 27074.875
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 27074.875
$ws.Range("N123").Value = -36874.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 40000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 40000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 40000
$ws.Range("N124").Value = -49820

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 35000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 35000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 35000
$ws.Range("N125").Value = -44840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 699.1667
$ws.Range("I126").Value = 671.8182
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2015.4546
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 454.5454
$ws.Range("N126").Value = -7940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 57000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 57000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 57000
$ws.Range("N127").Value = -66920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 49800
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 49800
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 49800
$ws.Range("N128").Value = -59760

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 49990
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49990
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49990
$ws.Range("N129").Value = -59990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 23525
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 23525
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 23525
$ws.Range("N130").Value = -33565

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 30000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 30000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1233.0886
$ws.Range("I132").Value = 755.2653
$ws.Range("J132").Value = 2013.5333
$ws.Range("K132").Value = 2265.7959
$ws.Range("L132").Value = 6040.5999
$ws.Range("M132").Value = 264.2040999999999
$ws.Range("N132").Value = -11100.5999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 76223
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 76223
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 76223
$ws.Range("N133").Value = -86343

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 112505
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 112505
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 112505
$ws.Range("N135").Value = -122645

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1553.2222
$ws.Range("I136").Value = 1526.3064
$ws.Range("J136").Value = 1641.0526
$ws.Range("K136").Value = 4578.9192
$ws.Range("L136").Value = 4923.1578
$ws.Range("M136").Value = -2028.9192
$ws.Range("N136").Value = -10023.1578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 60000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 60000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 75112.42999999999
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 75112.42999999999
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 75112.42999999999
$ws.Range("N138").Value = -85392.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 35496.727
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 35496.727
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 35496.727
$ws.Range("N140").Value = -45856.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 84590
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 84590
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 84590
$ws.Range("N141").Value = -94950
